$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Final, alphabetically-sorted list of test-game names that belongs in column
# A (row 1 is the header "9 Queen Checkmate", rows 2-76 are the sorted data).
# ---------------------------------------------------------------------------
$colA = @(
  '9 Queen Checkmate',
  '9 Queens',
  'AAA',
  'AQ',
  'Black Wins with Queen and EvilMorty',
  'Brian the dog',
  'Checkmate by promotion to Octocat',
  'Future Military tech',
  'Future war machines',
  'I',
  'II',
  'III',
  'IV',
  'NewCouncil',
  'NewCustom',
  'No pawn promo choices',
  'Stalemate 1',
  'Standard white wins using pawns',
  'Timberwolf and X Wing',
  'VII',
  'White about to win by capture black queen and promote',
  'White about to win with Queen in custom game',
  'White resigned the third',
  'White resigns',
  'White resigns again',
  'White wins 2 queen trap',
  'White wins by EvilMorty promotion',
  'White wins by dimond piece promotion',
  'White wins by pawn promo to rook in standard game',
  'White wins with 3 queens',
  'White wins with Evil Morty and Jester',
  'White wins with Evil Morty and Tux',
  'White wins with Queen Octocat Tux',
  'castle_test1',
  'castle_test2',
  'castle_test3',
  'castle_test4',
  'check_example1',
  'check_example2',
  'check_example3',
  'check_example4',
  'check_example6',
  'custom stalemate 1 black penguin trap white king',
  'dummy_game',
  'en-passant-test',
  'enpassant_test1',
  'fundemental_defense',
  'game with linux penguins',
  'king_range_test',
  'knight_threat1',
  'knight_threat2',
  'new_test',
  'pawn_hist_test',
  'pawn_moves_test',
  'pawn_promo',
  'pawn_range',
  'pawn_threat1',
  'pinned_ex1',
  'pinned_ex2',
  'promo_test',
  'ranges_test',
  'res1',
  'resignation',
  'stalemate only the Kings remain',
  'stalemate_exp1',
  'standard black pawn just en-passant',
  'standard stalemate 1 white bishops trap king',
  'super_checkmate_impossible_example',
  'test name save',
  'testing pawn hist again',
  'testing saving name of game',
  'two_queens',
  'undefined',
  'white resigns 2',
  'white wins by capturing black Queen and promoting to Queen',
  'white wins by dimond promotion again'
)

# Write the new column-A content (rows 1-76), overwriting whatever used to
# live there.
for ($i = 0; $i -lt $colA.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
}

# Row 73 also gets a "?" note in column B (new addition).
$ws.Cells.Item(73, 2).Value = "?"

# The sheet used to hold 111 rows of data; only 76 remain, so drop the tail
# and shift the remaining cells up so the used range shrinks to A1:B76.
$ws.Range("A77:A111").Delete(-4162)

# ---------------------------------------------------------------------------
# Column A formatting: keep the 67-character custom width but drop the
# column-level default style (individual cells keep their own style).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ClearFormats()
$ws.Columns.Item(1).ColumnWidth = 66.17
$ws.Range("A1:A76").Style = "Bad"

# Re-color the (builtin) "Bad" cell style so it reads as "Good" (green)
# instead of red - same cell style id, new colors.
$style = $wb.Styles.Item("Bad")
$style.Font.Color = 24832
$style.Interior.Color = 13561798

# Selection, matching the saved workbook state.
$ws.Range("E80").Select()

Write-Output "done"
